$wb = $excel.ActiveWorkbook

# OFF sheet ("OFF") - row 3 update
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 277
$wsOff.Range("C3").Value = 203
$wsOff.Range("D3").Value = 67
$wsOff.Range("E3").Value = 45
$wsOff.Range("F3").Value = 6
$wsOff.Range("G3").Value = 5

# DEF sheet ("DEF") - row 3 update
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 446
$wsDef.Range("C3").Value = 322
$wsDef.Range("D3").Value = 106
$wsDef.Range("E3").Value = 52
$wsDef.Range("F3").Value = 5
$wsDef.Range("G3").Value = 4
